$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 38 from 2023-10-22 (45221) to 2023-10-25 (45224)
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 3).Value = 45224
}
